$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C (rows 2-61) holds a date serial number that was bumped by one day
# (45203 -> 45204) in this automatic update.
for ($row = 2; $row -le 61; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value2 = $cell.Value2 + 1
}
